$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Pages" headings get their route annotated.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Страница регистрации:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Страница регистрации(/register):", 2) | Out-Null

$d.Content.Find.Execute("Страница авторизации:", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Страница авторизации(/authorize):", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) The "Рецепты:" bullet list gains two new fields: "id категории" and
#    "картинка". The existing "id автора рецепта" bullet is kept (it gets
#    duplicated as its own bullet right after the renamed one), matching the
#    diff which turns the old "id" + "автора рецепта" run pair into:
#       id  | категории          (renamed existing paragraph)
#       id  | автора рецепта     (new paragraph)
#       картинка                 (new paragraph)
# ---------------------------------------------------------------------------

# Helper: replace the text inside $range with $newText while forcing Word to
# keep it in its own run (instead of silently merging it into a
# neighbouring run that happens to carry identical formatting). This is done
# by toggling Bold on right before the assignment and back off right after.
function Set-RangeTextIsolated($range, $newText) {
    $rangeStart = $range.Start
    $range.Font.Bold = 1
    $range.Text = $newText
    $fixed = $d.Range($rangeStart, $rangeStart + $newText.Length)
    $fixed.Font.Bold = 0
    return $fixed
}

# Locate the paragraph whose whole text is "id автора рецепта" (Word
# terminates paragraph Range.Text with a trailing paragraph-mark char, so
# trim before comparing to stay robust to that).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r", "`n", "`a")
    if ($ptext -eq "id автора рецепта") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the 'id автора рецепта' bullet paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# -- 2a. Rename the second run of the existing bullet: "автора рецепта" -> "категории"
$idRunLength = 3   # "id " including the trailing space
$secondRunStart = $target.Range.Start + $idRunLength
$secondRunEnd = $target.Range.End - 1   # exclude the paragraph mark
$secondRun = $d.Range($secondRunStart, $secondRunEnd)
Set-RangeTextIsolated $secondRun "категории" | Out-Null

# Re-fetch the paragraph (indices/ranges stay valid here, but be defensive).
$target = $d.Paragraphs.Item($targetIndex)

# -- 2b. Insert a brand-new bullet paragraph right after it, re-creating the
#        original "id " + "автора рецепта" pair of runs.
$target.Range.InsertParagraphAfter() | Out-Null
$authorPara = $d.Paragraphs.Item($targetIndex + 1)
$authorPara.Range.Text = "id автора рецепта"
$authorPara = $d.Paragraphs.Item($targetIndex + 1)

$authorSecondStart = $authorPara.Range.Start + $idRunLength
$authorSecondEnd = $authorPara.Range.End - 1
$authorSecondRun = $d.Range($authorSecondStart, $authorSecondEnd)
Set-RangeTextIsolated $authorSecondRun "автора рецепта" | Out-Null

# -- 2c. Insert another new bullet paragraph after that one containing just "картинка".
$authorPara = $d.Paragraphs.Item($targetIndex + 1)
$authorPara.Range.InsertParagraphAfter() | Out-Null
$picturePara = $d.Paragraphs.Item($targetIndex + 2)
$picturePara.Range.Text = "картинка"

Write-Output "done"
